# FSR_S4 Stability Error Data - "did more analysis on fsrs4 day2"
#
# 1. Convert the H column (Max Error (lbf)) formulas for rows 2-22 into a
#    single shared formula group (mirrors the existing shared-formula
#    pattern already used by columns G and I).
# 2. Add a new data row (row 24) with a label for a second-day 5.00 run.
# 3. Update the view: zoom to 130% and move the active selection.
# 4. Widen column A slightly (it grew to fit the new, longer label).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-enter the H2:H22 formulas as one range write so the engine
#        collapses them into a shared formula group (t="shared"), just
#        like columns G and I already are. The formula text and the
#        calculated results are unchanged.
$ws.Range("H2:H22").Formula = "=A2*(1+(C2/100))"

# --- 2. New row: second-day 5.00 lbf stability point label.
$ws.Range("A24").Value = "5.00(2ndDay)"

# --- 3. View changes: zoom level and active cell/selection.
$excel.ActiveWindow.Zoom = 130
$ws.Range("G27").Select()

# --- 4. Column A grows from ~10.4 to ~12.6 chars to fit "5.00(2ndDay)".
$ws.Columns.Item(1).ColumnWidth = 11.6

$wb.Save()
